$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.030861739642096
$ws.Cells.Item(2, 4).Value = 1.035023953272568
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.02945350125632
$ws.Cells.Item(2, 9).Value = 1.035029052895697
$ws.Cells.Item(2, 10).Value = 1.036000645282986
$ws.Cells.Item(2, 11).Value = 1.037821745724467
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.032267367746458
$ws.Cells.Item(2, 14).Value = 1.015951985150247
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.031903166962125
$ws.Cells.Item(3, 4).Value = 1.035816745775685
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.031127590581169
$ws.Cells.Item(3, 9).Value = 1.035296101574965
$ws.Cells.Item(3, 10).Value = 1.036683236295808
$ws.Cells.Item(3, 11).Value = 1.03842389133978
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.03374727638474
$ws.Cells.Item(3, 14).Value = 1.016181718337603
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.032576334061871
$ws.Cells.Item(4, 4).Value = 1.036329107520588
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.032210033316676
$ws.Cells.Item(4, 9).Value = 1.035467392929081
$ws.Cells.Item(4, 10).Value = 1.037123653241846
$ws.Cells.Item(4, 11).Value = 1.038812240224125
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.034703620387131
$ws.Cells.Item(4, 14).Value = 1.016329858737484
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.032859166435399
$ws.Cells.Item(5, 4).Value = 1.036544354900518
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.032664907506443
$ws.Cells.Item(5, 9).Value = 1.03553904360797
$ws.Cells.Item(5, 10).Value = 1.037308503078192
$ws.Cells.Item(5, 11).Value = 1.038975196801805
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.035105373752261
$ws.Cells.Item(5, 14).Value = 1.016392014702252
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.032906645494394
$ws.Cells.Item(6, 4).Value = 1.036580487140617
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.032741272282796
$ws.Cells.Item(6, 9).Value = 1.035551052966084
$ws.Cells.Item(6, 10).Value = 1.037339522542972
$ws.Cells.Item(6, 11).Value = 1.039002540055828
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.03517281289291
$ws.Cells.Item(6, 14).Value = 1.01640244380427
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.032580113935969
$ws.Cells.Item(7, 4).Value = 1.036331984253001
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.032216112083858
$ws.Cells.Item(7, 9).Value = 1.035468351743621
$ws.Cells.Item(7, 10).Value = 1.037126124397147
$ws.Cells.Item(7, 11).Value = 1.038814418856205
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.034708989780453
$ws.Cells.Item(7, 14).Value = 1.016330689748438
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.031213842241225
$ws.Cells.Item(8, 4).Value = 1.035292011762318
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.0300194390763
$ws.Cells.Item(8, 9).Value = 1.035119615627194
$ws.Cells.Item(8, 10).Value = 1.036231592911139
$ws.Cells.Item(8, 11).Value = 1.038025509234342
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.03276777405136
$ws.Cells.Item(8, 14).Value = 1.016029730877587
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.028800801700606
$ws.Cells.Item(9, 4).Value = 1.033454598360479
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.026142097084033
$ws.Cells.Item(9, 9).Value = 1.034493530165216
$ws.Cells.Item(9, 10).Value = 1.034645560932856
$ws.Cells.Item(9, 11).Value = 1.036625501984826
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.029337171548478
$ws.Cells.Item(9, 14).Value = 1.015495461030824
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.02718828727175
$ws.Cells.Item(10, 4).Value = 1.032226338036214
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.023552314334695
$ws.Cells.Item(10, 9).Value = 1.034068324411653
$ws.Cells.Item(10, 10).Value = 1.033581556598046
$ws.Cells.Item(10, 11).Value = 1.03568547386455
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.027042984623606
$ws.Cells.Item(10, 14).Value = 1.015136604706191
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.026489115531423
$ws.Cells.Item(11, 4).Value = 1.031693686913018
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.022429638328499
$ws.Cells.Item(11, 9).Value = 1.033882343025598
$ws.Cells.Item(11, 10).Value = 1.033119232491233
$ws.Cells.Item(11, 11).Value = 1.035276828447193
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.026047789439016
$ws.Cells.Item(11, 14).Value = 1.0149805750223
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.026229267834599
$ws.Cells.Item(12, 4).Value = 1.031495714132627
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.02201242400102
$ws.Cells.Item(12, 9).Value = 1.033812980309371
$ws.Cells.Item(12, 10).Value = 1.032947261942295
$ws.Cells.Item(12, 11).Value = 1.035124796424328
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.025677851087005
$ws.Cells.Item(12, 14).Value = 1.014922521553494
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.026285012614216
$ws.Cells.Item(13, 4).Value = 1.031538185527658
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.022101927209492
$ws.Cells.Item(13, 9).Value = 1.033827871574394
$ws.Cells.Item(13, 10).Value = 1.03298416122673
$ws.Cells.Item(13, 11).Value = 1.035157418825901
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.025757216862181
$ws.Cells.Item(13, 14).Value = 1.014934978622599
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.026467639396148
$ws.Cells.Item(14, 4).Value = 1.031677324926496
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.022395155448834
$ws.Cells.Item(14, 9).Value = 1.033876615217935
$ws.Cells.Item(14, 10).Value = 1.033105022316078
$ws.Cells.Item(14, 11).Value = 1.035264266396871
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.02601721595144
$ws.Cells.Item(14, 14).Value = 1.014975778288496
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.026580142600183
$ws.Cells.Item(15, 4).Value = 1.031763037039006
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.022575795923615
$ws.Cells.Item(15, 9).Value = 1.033906610532069
$ws.Cells.Item(15, 10).Value = 1.033179456592071
$ws.Cells.Item(15, 11).Value = 1.035330066477942
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.026177372804031
$ws.Cells.Item(15, 14).Value = 1.01500090342592
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.027234668921427
$ws.Cells.Item(16, 4).Value = 1.032261671238739
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.023626794784731
$ws.Cells.Item(16, 9).Value = 1.034080628030412
$ws.Cells.Item(16, 10).Value = 1.033612205584587
$ws.Cells.Item(16, 11).Value = 1.035712560330607
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.027108993858534
$ws.Cells.Item(16, 14).Value = 1.015146946297096
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027644981865159
$ws.Cells.Item(17, 4).Value = 1.032574234518626
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.024285708273516
$ws.Cells.Item(17, 9).Value = 1.03418928484978
$ws.Cells.Item(17, 10).Value = 1.033883227016275
$ws.Cells.Item(17, 11).Value = 1.035952057329716
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.027692887976928
$ws.Cells.Item(17, 14).Value = 1.015238382712772
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.027884219715903
$ws.Cells.Item(18, 4).Value = 1.032756469683875
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.024669918229219
$ws.Cells.Item(18, 9).Value = 1.034252482694195
$ws.Cells.Item(18, 10).Value = 1.034041154675939
$ws.Cells.Item(18, 11).Value = 1.03609159692168
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.028033290719466
$ws.Cells.Item(18, 14).Value = 1.015291654070737
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.0279657783271
$ws.Cells.Item(19, 4).Value = 1.03281859402451
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.024800903304189
$ws.Cells.Item(19, 9).Value = 1.034274001047112
$ws.Cells.Item(19, 10).Value = 1.034094977769601
$ws.Cells.Item(19, 11).Value = 1.036139150078138
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.028149330204351
$ws.Cells.Item(19, 14).Value = 1.01530980773535
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.027600968574318
$ws.Cells.Item(20, 4).Value = 1.032540707483564
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.024215025897543
$ws.Cells.Item(20, 9).Value = 1.034177645608928
$ws.Cells.Item(20, 10).Value = 1.033854164981658
$ws.Cells.Item(20, 11).Value = 1.035926377606518
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.027630259600145
$ws.Cells.Item(20, 14).Value = 1.015228578863128
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.026413864361866
$ws.Cells.Item(21, 4).Value = 1.031636355231889
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.022308812708431
$ws.Cells.Item(21, 9).Value = 1.033862269192081
$ws.Cells.Item(21, 10).Value = 1.033069438455424
$ws.Cells.Item(21, 11).Value = 1.035232809168076
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.025940660458263
$ws.Cells.Item(21, 14).Value = 1.014963766488714
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.02566664933088
$ws.Cells.Item(22, 4).Value = 1.031067044037728
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.021109124523065
$ws.Cells.Item(22, 9).Value = 1.033662353622409
$ws.Cells.Item(22, 10).Value = 1.032574644063289
$ws.Cells.Item(22, 11).Value = 1.034795328958426
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.024876726540598
$ws.Cells.Item(22, 14).Value = 1.014796706380914
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.026062842227507
$ws.Cells.Item(23, 4).Value = 1.031368914321541
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.0217452166908
$ws.Cells.Item(23, 9).Value = 1.033768487044385
$ws.Cells.Item(23, 10).Value = 1.032837077823001
$ws.Cells.Item(23, 11).Value = 1.035027379285845
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.025440894079257
$ws.Cells.Item(23, 14).Value = 1.014885321571867
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.027620856565632
$ws.Cells.Item(24, 4).Value = 1.03255585714725
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.024246964598998
$ws.Cells.Item(24, 9).Value = 1.034182905436294
$ws.Cells.Item(24, 10).Value = 1.033867297341168
$ws.Cells.Item(24, 11).Value = 1.03593798164642
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.027658559199801
$ws.Cells.Item(24, 14).Value = 1.015233008992053
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.029425294645133
$ws.Cells.Item(25, 4).Value = 1.033930194291949
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.027145307509501
$ws.Cells.Item(25, 9).Value = 1.034656762987444
$ws.Cells.Item(25, 10).Value = 1.0350567534524
$ws.Cells.Item(25, 11).Value = 1.036988611461288
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.030225286770668
$ws.Cells.Item(25, 14).Value = 1.01563405239268
